$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Rarres2"
$ws.Cells.Item(2,3).Value = "Gpr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.7062263333333333
$ws.Cells.Item(2,8).Value = 2.118679
$ws.Cells.Item(2,9).Value = 0.03412450835304152
$ws.Cells.Item(2,10).Value = 0.03412450835304151
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.576048333333333
$ws.Cells.Item(2,14).Value = 4.728145
$ws.Cells.Item(2,15).Value = 0.3961230229224872
$ws.Cells.Item(2,16).Value = 0.3961230229224872
$ws.Cells.Item(2,17).Value = 1.113046835606111
$ws.Cells.Item(2,18).Value = 10.017421520455
$ws.Cells.Item(2,19).Value = 0.01351750340455047
$ws.Cells.Item(2,20).Value = 0.01351750340455047

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Rarres2"
$ws.Cells.Item(3,3).Value = "Gpr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.7062263333333333
$ws.Cells.Item(3,8).Value = 2.118679
$ws.Cells.Item(3,9).Value = 0.03412450835304152
$ws.Cells.Item(3,10).Value = 0.03412450835304151
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.631411
$ws.Cells.Item(3,14).Value = 4.894233
$ws.Cells.Item(3,15).Value = 0.4100378416581965
$ws.Cells.Item(3,16).Value = 0.4100378416581965
$ws.Cells.Item(3,17).Value = 1.152145408689667
$ws.Cells.Item(3,18).Value = 10.369308678207
$ws.Cells.Item(3,19).Value = 0.01399233975272824
$ws.Cells.Item(3,20).Value = 0.01399233975272824

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Rarres2"
$ws.Cells.Item(4,3).Value = "Gpr1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.7062263333333333
$ws.Cells.Item(4,8).Value = 2.118679
$ws.Cells.Item(4,9).Value = 0.03412450835304152
$ws.Cells.Item(4,10).Value = 0.03412450835304151
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.5138943333333333
$ws.Cells.Item(4,14).Value = 1.541683
$ws.Cells.Item(4,15).Value = 0.1291618870293125
$ws.Cells.Item(4,16).Value = 0.1291618870293125
$ws.Cells.Item(4,17).Value = 0.3629257107507778
$ws.Cells.Item(4,18).Value = 3.266331396756999
$ws.Cells.Item(4,19).Value = 0.00440758589282638
$ws.Cells.Item(4,20).Value = 0.004407585892826379

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Rarres2"
$ws.Cells.Item(5,3).Value = "Gpr1"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.7062263333333333
$ws.Cells.Item(5,8).Value = 2.118679
$ws.Cells.Item(5,9).Value = 0.03412450835304152
$ws.Cells.Item(5,10).Value = 0.03412450835304151
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.2573303333333333
$ws.Cells.Item(5,14).Value = 0.771991
$ws.Cells.Item(5,15).Value = 0.06467724839000366
$ws.Cells.Item(5,16).Value = 0.06467724839000366
$ws.Cells.Item(5,17).Value = 0.1817334577654444
$ws.Cells.Item(5,18).Value = 1.635601119889
$ws.Cells.Item(5,19).Value = 0.002207079302936421
$ws.Cells.Item(5,20).Value = 0.002207079302936421

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Rarres2"
$ws.Cells.Item(6,3).Value = "Gpr1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 18.279338
$ws.Cells.Item(6,8).Value = 54.838014
$ws.Cells.Item(6,9).Value = 0.883248602930037
$ws.Cells.Item(6,10).Value = 0.8832486029300368
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 1.576048333333333
$ws.Cells.Item(6,14).Value = 4.728145
$ws.Cells.Item(6,15).Value = 0.3961230229224872
$ws.Cells.Item(6,16).Value = 0.3961230229224872
$ws.Cells.Item(6,17).Value = 28.80912018933666
$ws.Cells.Item(6,18).Value = 259.28208170403
$ws.Cells.Item(6,19).Value = 0.3498751065847099
$ws.Cells.Item(6,20).Value = 0.3498751065847098

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Rarres2"
$ws.Cells.Item(7,3).Value = "Gpr1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 18.279338
$ws.Cells.Item(7,8).Value = 54.838014
$ws.Cells.Item(7,9).Value = 0.883248602930037
$ws.Cells.Item(7,10).Value = 0.8832486029300368
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.631411
$ws.Cells.Item(7,14).Value = 4.894233
$ws.Cells.Item(7,15).Value = 0.4100378416581965
$ws.Cells.Item(7,16).Value = 0.4100378416581965
$ws.Cells.Item(7,17).Value = 29.821113085918
$ws.Cells.Item(7,18).Value = 268.390017773262
$ws.Cells.Item(7,19).Value = 0.3621653507930498
$ws.Cells.Item(7,20).Value = 0.3621653507930497

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Rarres2"
$ws.Cells.Item(8,3).Value = "Gpr1"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 18.279338
$ws.Cells.Item(8,8).Value = 54.838014
$ws.Cells.Item(8,9).Value = 0.883248602930037
$ws.Cells.Item(8,10).Value = 0.8832486029300368
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.5138943333333333
$ws.Cells.Item(8,14).Value = 1.541683
$ws.Cells.Item(8,15).Value = 0.1291618870293125
$ws.Cells.Item(8,16).Value = 0.1291618870293125
$ws.Cells.Item(8,17).Value = 9.393648215284667
$ws.Cells.Item(8,18).Value = 84.542833937562
$ws.Cells.Item(8,19).Value = 0.1140820562704476
$ws.Cells.Item(8,20).Value = 0.1140820562704475

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Rarres2"
$ws.Cells.Item(9,3).Value = "Gpr1"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 18.279338
$ws.Cells.Item(9,8).Value = 54.838014
$ws.Cells.Item(9,9).Value = 0.883248602930037
$ws.Cells.Item(9,10).Value = 0.8832486029300368
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.2573303333333333
$ws.Cells.Item(9,14).Value = 0.771991
$ws.Cells.Item(9,15).Value = 0.06467724839000366
$ws.Cells.Item(9,16).Value = 0.06467724839000366
$ws.Cells.Item(9,17).Value = 4.703828140652666
$ws.Cells.Item(9,18).Value = 42.334453265874
$ws.Cells.Item(9,19).Value = 0.05712608928182972
$ws.Cells.Item(9,20).Value = 0.0571260892818297

# Row 10
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Rarres2"
$ws.Cells.Item(10,3).Value = "Gpr1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.03164166666666667
$ws.Cells.Item(10,8).Value = 0.094925
$ws.Cells.Item(10,9).Value = 0.001528909738290919
$ws.Cells.Item(10,10).Value = 0.001528909738290919
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 1.576048333333333
$ws.Cells.Item(10,14).Value = 4.728145
$ws.Cells.Item(10,15).Value = 0.3961230229224872
$ws.Cells.Item(10,16).Value = 0.3961230229224872
$ws.Cells.Item(10,17).Value = 0.04986879601388888
$ws.Cells.Item(10,18).Value = 0.4488191641249999
$ws.Cells.Item(10,19).Value = 0.0006056363473074278
$ws.Cells.Item(10,20).Value = 0.0006056363473074276

# Row 11
$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Rarres2"
$ws.Cells.Item(11,3).Value = "Gpr1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.03164166666666667
$ws.Cells.Item(11,8).Value = 0.094925
$ws.Cells.Item(11,9).Value = 0.001528909738290919
$ws.Cells.Item(11,10).Value = 0.001528909738290919
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.631411
$ws.Cells.Item(11,14).Value = 4.894233
$ws.Cells.Item(11,15).Value = 0.4100378416581965
$ws.Cells.Item(11,16).Value = 0.4100378416581965
$ws.Cells.Item(11,17).Value = 0.05162056305833333
$ws.Cells.Item(11,18).Value = 0.464585067525
$ws.Cells.Item(11,19).Value = 0.0006269108491790066
$ws.Cells.Item(11,20).Value = 0.0006269108491790065

# Row 12
$ws.Cells.Item(12,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,2).Value = "Rarres2"
$ws.Cells.Item(12,3).Value = "Gpr1"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.03164166666666667
$ws.Cells.Item(12,8).Value = 0.094925
$ws.Cells.Item(12,9).Value = 0.001528909738290919
$ws.Cells.Item(12,10).Value = 0.001528909738290919
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.5138943333333333
$ws.Cells.Item(12,14).Value = 1.541683
$ws.Cells.Item(12,15).Value = 0.1291618870293125
$ws.Cells.Item(12,16).Value = 0.1291618870293125
$ws.Cells.Item(12,17).Value = 0.01626047319722222
$ws.Cells.Item(12,18).Value = 0.146344258775
$ws.Cells.Item(12,19).Value = 0.0001974768668951475
$ws.Cells.Item(12,20).Value = 0.0001974768668951474

# Row 13
$ws.Cells.Item(13,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13,2).Value = "Rarres2"
$ws.Cells.Item(13,3).Value = "Gpr1"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.03164166666666667
$ws.Cells.Item(13,8).Value = 0.094925
$ws.Cells.Item(13,9).Value = 0.001528909738290919
$ws.Cells.Item(13,10).Value = 0.001528909738290919
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.2573303333333333
$ws.Cells.Item(13,14).Value = 0.771991
$ws.Cells.Item(13,15).Value = 0.06467724839000366
$ws.Cells.Item(13,16).Value = 0.06467724839000366
$ws.Cells.Item(13,17).Value = 0.008142360630555555
$ws.Cells.Item(13,18).Value = 0.073281245675
$ws.Cells.Item(13,19).Value = 0.00009888567490933728
$ws.Cells.Item(13,20).Value = 0.00009888567490933725

# Row 14
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Rarres2"
$ws.Cells.Item(14,3).Value = "Gpr1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 1.678369333333333
$ws.Cells.Item(14,8).Value = 5.035108
$ws.Cells.Item(14,9).Value = 0.08109797897863065
$ws.Cells.Item(14,10).Value = 0.08109797897863062
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 1.576048333333333
$ws.Cells.Item(14,14).Value = 4.728145
$ws.Cells.Item(14,15).Value = 0.3961230229224872
$ws.Cells.Item(14,16).Value = 0.3961230229224872
$ws.Cells.Item(14,17).Value = 2.645191190517778
$ws.Cells.Item(14,18).Value = 23.80672071466
$ws.Cells.Item(14,19).Value = 0.03212477658591949
$ws.Cells.Item(14,20).Value = 0.03212477658591949

# Row 15
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Rarres2"
$ws.Cells.Item(15,3).Value = "Gpr1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 1.678369333333333
$ws.Cells.Item(15,8).Value = 5.035108
$ws.Cells.Item(15,9).Value = 0.08109797897863065
$ws.Cells.Item(15,10).Value = 0.08109797897863062
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 1.631411
$ws.Cells.Item(15,14).Value = 4.894233
$ws.Cells.Item(15,15).Value = 0.4100378416581965
$ws.Cells.Item(15,16).Value = 0.4100378416581965
$ws.Cells.Item(15,17).Value = 2.738110192462667
$ws.Cells.Item(15,18).Value = 24.642991732164
$ws.Cells.Item(15,19).Value = 0.0332532402632395
$ws.Cells.Item(15,20).Value = 0.03325324026323949

# Row 16
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Rarres2"
$ws.Cells.Item(16,3).Value = "Gpr1"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 1.678369333333333
$ws.Cells.Item(16,8).Value = 5.035108
$ws.Cells.Item(16,9).Value = 0.08109797897863065
$ws.Cells.Item(16,10).Value = 0.08109797897863062
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.5138943333333333
$ws.Cells.Item(16,14).Value = 1.541683
$ws.Cells.Item(16,15).Value = 0.1291618870293125
$ws.Cells.Item(16,16).Value = 0.1291618870293125
$ws.Cells.Item(16,17).Value = 0.8625044896404445
$ws.Cells.Item(16,18).Value = 7.762540406764
$ws.Cells.Item(16,19).Value = 0.01047476799914345
$ws.Cells.Item(16,20).Value = 0.01047476799914345

# Row 17
$ws.Cells.Item(17,1).Value = "MuSCs"
$ws.Cells.Item(17,2).Value = "Rarres2"
$ws.Cells.Item(17,3).Value = "Gpr1"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 1.678369333333333
$ws.Cells.Item(17,8).Value = 5.035108
$ws.Cells.Item(17,9).Value = 0.08109797897863065
$ws.Cells.Item(17,10).Value = 0.08109797897863062
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.2573303333333333
$ws.Cells.Item(17,14).Value = 0.771991
$ws.Cells.Item(17,15).Value = 0.06467724839000366
$ws.Cells.Item(17,16).Value = 0.06467724839000366
$ws.Cells.Item(17,17).Value = 0.4318953400031111
$ws.Cells.Item(17,18).Value = 3.887058060028
$ws.Cells.Item(17,19).Value = 0.00524519413032819
$ws.Cells.Item(17,20).Value = 0.005245194130328188

Write-Host "done"